# Updated legacy GSC export data: drop the oldest day (2025-10-02) from the
# rolling coverage export on the "Chart" sheet so every later day's row
# shifts up by one. (Critical issues / Non-critical issues / Metadata
# sheets are untouched; only the shared-string table shrinks as a side
# effect of Excel's normal bookkeeping.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows(2).Delete()
